{"js": "// Apply the three grammatical / wording fixes described by the diff:\n//  1. \" (among ATO, AFP, Border control and foreign counterparts). \"\n//       -> \" among ATO, AFP, Immigration and AUSTRAC. \"\n//  2. \"have to looking for clues\"  -> \"have to look for clues\"\n//  3. \"digitalizing conventional audition practice\"\n//       -> \"digitalizing conventional tax audit practice\"\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replacementText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacementText, \"Replace\");\n  await context.sync();\n}\n\nawait replaceOnce(\n  \" (among ATO, AFP, Border control and foreign counterparts). \",\n  \" among ATO, AFP, Immigration and AUSTRAC. \"\n);\n\nawait replaceOnce(\n  \"have to looking for clues\",\n  \"have to look for clues\"\n);\n\nawait replaceOnce(\n  \"digitalizing conventional audition practice\",\n  \"digitalizing conventional tax audit practice\"\n);\n", "ps1": "# Apply the three grammatical / wording fixes described by the diff:\n#  1. \" (among ATO, AFP, Border control and foreign counterparts). \"\n#       -> \" among ATO, AFP, Immigration and AUSTRAC. \"\n#  2. \"have to looking for clues\"  -> \"have to look for clues\"\n#  3. \"digitalizing conventional audition practice\"\n#       -> \"digitalizing conventional tax audit practice\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, $find.Replacement.Text, 2) | Out-Null  # wdReplaceAll\n}\n\nReplace-Text \" (among ATO, AFP, Border control and foreign counterparts). \" \" among ATO, AFP, Immigration and AUSTRAC. \"\nReplace-Text \"have to looking for clues\" \"have to look for clues\"\nReplace-Text \"digitalizing conventional audition practice\" \"digitalizing conventional tax audit practice\"\n"}
